# Refreshed regression predictions (age 19-94) and dropped the trailing
# age=94-only duplicate row that the old plot data carried (age=18 dropped,
# new max age=94 now ends one row earlier). Also re-fits the underlying
# regression coefficients (ages now treated as a factor), so every wage
# prediction in column B is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ages  = @(19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94)
$wages = @(4740.582554633078, 4860.25429175882, 4977.447166116676, 5091.839586717203, 5203.111772889624, 5310.947223449562, 5415.03420049433, 5515.067218218306, 5610.748526897786, 5701.789582033169, 5787.912488562154, 5868.851410071082, 5944.353933034749, 6014.182376309222, 6078.115036385573, 6135.947359284442, 6187.493030428849, 6232.584974372676, 6271.076256879731, 6302.840882538548, 6327.774481853854, 6345.794882571556, 6356.842560860144, 6360.880968881814, 6357.89673623036, 6347.899743682322, 6330.923068691606, 6307.022803048764, 6276.277744111486, 6238.788961984935, 6194.67924597918, 6144.092434586227, 6087.192634093319, 6024.163331772217, 5955.206410350143, 5880.541071168359, 5800.402674062861, 5715.041502553963, 5624.721463401599, 5529.718729967993, 5430.32033912655, 5326.822751663119, 5219.530386232591, 5108.754136961149, 4994.809884722841, 4878.017011970586, 4758.696930770997, 4637.171633381299, 4513.762274321521, 4388.787792440849, 4262.56358095994, 4135.400212897344, 4007.602228665517, 3879.46699195734, 3751.28361934552, 3623.331988292178, 3495.881827522673, 3369.191892963853, 3243.509231690026, 3119.068535567904, 2996.091585551384, 2874.786786855084, 2755.348794538934, 2637.958228369727, 2522.781475195773, 2409.970576481152, 2299.663198101779, 2191.982679009046, 2087.038154921419, 1984.924752811954, 1885.723851621767, 1789.5034043467, 1696.318316416809, 1606.210875115848, 1519.21122466882, 1435.337881559208)

# Rows 2..77 get the new age / predicted-wage pairs (overwrites the old
# rows 2..77 in place, so ages shift from 18-93 to 19-94).
for ($i = 0; $i -lt $ages.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $ages[$i]
    $ws.Cells.Item($r, 2).Value = $wages[$i]
}

# The data set is now one row shorter (78 -> 77 data+header rows), so the
# former last row (age=94 / old trailing wage) is removed entirely.
$ws.Range("A78:B78").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
